# Update Name of Algo
# Apply updated imputed values (columns A and C) for the RandomForest result sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "A4"   = -21.44439999999999
    "A6"   = -22.6513
    "A7"   = -19.39129999999999
    "C7"   = -12.2014
    "C12"  = -10.6587
    "C15"  = -14.8203
    "A16"  = -21.50689999999999
    "A20"  = -19.1984
    "C20"  = -11.9602
    "C21"  = -11.9834
    "C22"  = -12.54820000000001
    "C23"  = -12.2061
    "A28"  = -22.11429999999999
    "A29"  = -20.87089999999997
    "C29"  = -11.49950000000001
    "A32"  = -21.1771
    "C34"  = -11.25450000000002
    "A40"  = -20.246
    "C42"  = -11.91879999999999
    "C43"  = -12.75799999999999
    "C44"  = -13.76059999999999
    "C45"  = -13.42859999999999
    "A46"  = -21.9009
    "C46"  = -13.31889999999999
    "C50"  = -13.78549999999998
    "A51"  = -22.13090000000001
    "C51"  = -12.42209999999999
    "A52"  = -22.0907
    "A57"  = -22.68760000000001
    "A59"  = -22.0645
    "A62"  = -22.0763
    "A66"  = -22.1368
    "C66"  = -12.79769999999999
    "C67"  = -11.0591
    "A73"  = -20.81280000000002
    "A74"  = -21.84939999999998
    "C79"  = -11.16680000000001
    "C84"  = -13.61829999999998
    "A92"  = -21.80359999999999
    "C92"  = -11.2756
    "C97"  = -11.50060000000001
    "A100" = -22.2018
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
